# "Add files via upload" — replaces the two rich-value (linked image)
# cells C2/C3, which evaluated to #VALUE! errors, with plain URL text,
# fills in the previously-empty A4/B4 text cells, and moves the active
# selection from D3 to C4.
#
# Shared-string insertion order matters for an exact OOXML match, so we
# write A4/B4 before C2/C3 (matches the target sharedStrings.xml order:
# ... "i eat halouf?", "nah i dont think so",
#     "https://plus.unsplash.com/...", "https://cdn.pixabay.com/...").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "i eat halouf?"
$ws.Range("B4").Value = "nah i dont think so"

$ws.Range("C2").Value = "https://plus.unsplash.com/premium_photo-1664474619075-644dd191935f?fm=jpg&q=60&w=3000&ixlib=rb-4.0.3&ixid=M3wxMjA3fDB8MHxzZWFyY2h8MXx8aW1hZ2V8ZW58MHx8MHx8fDA%3D"
$ws.Range("C3").Value = "https://cdn.pixabay.com/photo/2015/04/23/22/00/tree-736885_1280.jpg"

# Move the selection from D3 to C4, as in the saved workbook view state.
$ws.Range("C4").Select() | Out-Null
